# Weekly crime-stat refresh: bump the report week/volume header text and
# overwrite the crime-count table (rows 14-29) with the newly collected data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (Volume/Number banner + "Report Covering the Week" line) ---
$ws.Range("A8").Value  = "Volume 29   Number  49"
$ws.Range("C9").Value  = "Report Covering the Week  12/5/2022  Through  12/11/2022"

# Helper: write a numeric cell using the "#,##0" integer-count style (same
# look as the table's existing count columns).
function Set-Count($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $value
    $cell.NumberFormat = "#,##0"
}

# Helper: write a numeric cell using the "#,##0.0;""-""#,##0.0" style used
# for the %-change columns.
function Set-Pct($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $value
    $cell.NumberFormat = '#,##0.0;"-"#,##0.0'
}

# Helper: the bolded TOTAL row (21) uses a 2-decimal variant of the
# %-change style ("#,##0.00;""-""#,##0.00") - keep it distinct so the row
# keeps its existing (bold) TOTAL-row style instead of sliding onto the
# regular row style.
function Set-PctTotal($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $value
    $cell.NumberFormat = '#,##0.00;"-"#,##0.00'
}

# Column map: C=WTD-2022 D=WTD-2021 E=WTD-%chg F=28day-2022 G=28day-2021
#             H=28day-%chg I=YTD-2022 J=YTD-2021 K=YTD-%chg L=2yr-%chg
#             M=12yr-%chg  N=29yr-%chg

# Row 14 - Murder
Set-Count 14 3  1
Set-Count 14 6  1
Set-Count 14 9  5
Set-Pct   14 11 0
Set-Pct   14 12 66.666666666666
Set-Pct   14 13 25
Set-Pct   14 14 -76.190476190476

# Row 16 - Robbery
Set-Count 16 3  7
Set-Pct   16 5  40
Set-Count 16 6  15
Set-Count 16 7  18
Set-Pct   16 8  -16.666666666666
Set-Count 16 9  216
Set-Count 16 10 184
Set-Pct   16 11 17.391304347826
Set-Pct   16 12 27.810650887574
Set-Pct   16 13 -13.253012048192
Set-Pct   16 14 -74.069627851140

# Row 17 - Fel. Assault
Set-Count 17 3  3
Set-Pct   17 5  -40
Set-Count 17 7  24
Set-Pct   17 8  -29.166666666666
Set-Count 17 9  323
Set-Count 17 10 328
Set-Pct   17 11 -1.524390243902
Set-Pct   17 12 19.629629629629
Set-Pct   17 13 125.874125874126
Set-Pct   17 14 -18.020304568527

# Row 18 - Burglary
Set-Count 18 3  2
Set-Count 18 4  4
Set-Pct   18 5  -50
Set-Count 18 6  6
Set-Count 18 7  16
Set-Pct   18 8  -62.5
Set-Count 18 9  140
Set-Count 18 10 109
Set-Pct   18 11 28.440366972477
Set-Pct   18 12 -9.677419354838
Set-Pct   18 13 -51.388888888888
Set-Pct   18 14 -87.868284228769

# Row 19 - Gr. Larceny
Set-Count 19 3  20
Set-Count 19 4  8
Set-Pct   19 5  150
Set-Count 19 6  55
Set-Count 19 7  45
Set-Pct   19 8  22.222222222222
Set-Count 19 9  631
Set-Count 19 10 485
Set-Pct   19 11 30.103092783505
Set-Pct   19 12 25.198412698412
Set-Pct   19 13 82.369942196531
Set-Pct   19 14 12.477718360071

# Row 20 - G.L.A.
Set-Count 20 3  5
Set-Pct   20 5  0
Set-Count 20 6  17
Set-Pct   20 8  -43.333333333333
Set-Count 20 9  212
Set-Count 20 10 200
Set-Pct   20 11 6
Set-Pct   20 12 18.435754189944
Set-Pct   20 13 -30.944625407166
Set-Pct   20 14 -93.239795918367

# Row 21 - TOTAL
Set-Count    21 3  38
Set-Count    21 4  27
Set-PctTotal 21 5  40.740740740740
Set-Count    21 6  111
Set-Count    21 7  134
Set-PctTotal 21 8  -17.164179104477
Set-Count    21 9  1546
Set-Count    21 10 1323
Set-PctTotal 21 11 16.855631141345
Set-PctTotal 21 12 19.844961240310
Set-PctTotal 21 13 14.095940959409
Set-PctTotal 21 14 -74.779771615008

# Row 22 - Transit
Set-Count 22 4  2
Set-Pct   22 5  -100
Set-Count 22 6  1
Set-Count 22 7  2
Set-Pct   22 8  -50
Set-Count 22 10 16
Set-Pct   22 11 56.25

# Row 24 - Petit Larceny
Set-Count 24 3  29
Set-Count 24 4  15
Set-Pct   24 5  93.333333333333
Set-Count 24 6  89
Set-Count 24 7  98
Set-Pct   24 8  -9.183673469387
Set-Count 24 9  1293
Set-Count 24 10 923
Set-Pct   24 11 40.086673889490
Set-Pct   24 12 59.432799013563
Set-Pct   24 13 95.317220543806

# Row 25 - Misd. Assault
Set-Count 25 3  10
Set-Count 25 4  9
Set-Pct   25 5  11.111111111111
Set-Count 25 7  28
Set-Pct   25 8  42.857142857142
Set-Count 25 9  507
Set-Count 25 10 440
Set-Pct   25 11 15.227272727272
Set-Pct   25 12 20.142180094786
Set-Pct   25 13 7.188160676532

# Row 26 - UCR Rape*
Set-Count 26 6  1
Set-Pct   26 8  0

# Row 27 - Other Sex Crimes
Set-Count 27 3  1
Set-Pct   27 5  0
Set-Count 27 6  3
Set-Pct   27 8  -25
Set-Count 27 9  55
Set-Count 27 10 39
Set-Pct   27 11 41.025641025641
Set-Pct   27 12 89.655172413793

# Row 28 - Shooting Vic.
Set-Pct   28 12 128.571428571429

# Row 29 - Shooting Inc.
Set-Pct   29 12 100
